$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.196.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.789.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "116.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.131"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.227.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.792.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.888"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.075.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0982"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0825"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0409"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.86%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  +20.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.116"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.073.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.910"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
